# Refresh cryptos list: updated Price (D) and Volume(1h) (E) text cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Leading apostrophe forces literal-text entry so Excel doesn't
    # re-interpret numeric-looking strings (e.g. "75.00", "6.99") as
    # numbers; resetting the style back to Normal afterwards drops the
    # transient quote-prefix flag Excel stamps on the cell format.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '43.203.41'
Set-TextValue $ws.Range("E2") '  -1.35%  '
Set-TextValue $ws.Range("D3") '2.268.64'
Set-TextValue $ws.Range("E3") '  -1.79%  '
Set-TextValue $ws.Range("E4") '  +0.08%  '
Set-TextValue $ws.Range("D5") '113.38'
Set-TextValue $ws.Range("E5") '  +4.50%  '
Set-TextValue $ws.Range("D6") '264.74'
Set-TextValue $ws.Range("E6") '  -2.51%  '
Set-TextValue $ws.Range("E7") '  -1.19%  '
Set-TextValue $ws.Range("D9") '0.598'
Set-TextValue $ws.Range("E9") '  -3.38%  '
Set-TextValue $ws.Range("D10") '47.98'
Set-TextValue $ws.Range("E10") '  +0.20%  '
Set-TextValue $ws.Range("E11") '  -1.75%  '
Set-TextValue $ws.Range("D12") '8.74'
Set-TextValue $ws.Range("E12") '  +3.73%  '
Set-TextValue $ws.Range("E13") '  -0.40%  '
Set-TextValue $ws.Range("D14") '15.38'
Set-TextValue $ws.Range("E14") '  -2.80%  '
Set-TextValue $ws.Range("D15") '2.605.19'
Set-TextValue $ws.Range("E15") '  -1.62%  '
Set-TextValue $ws.Range("E16") '  -0.80%  '
Set-TextValue $ws.Range("D17") '2.277.80'
Set-TextValue $ws.Range("E17") '  -1.19%  '
Set-TextValue $ws.Range("D18") '43.130.22'
Set-TextValue $ws.Range("E18") '  -1.49%  '
Set-TextValue $ws.Range("E19") '  -3.61%  '
Set-TextValue $ws.Range("D20") '6.99'
Set-TextValue $ws.Range("E20") '  +10.57%  '
Set-TextValue $ws.Range("D21") '71.13'
Set-TextValue $ws.Range("E21") '  -1.72%  '
Set-TextValue $ws.Range("E22") '  -3.31%  '
Set-TextValue $ws.Range("D23") '9.83'
Set-TextValue $ws.Range("E23") '  +6.25%  '
Set-TextValue $ws.Range("D24") '230.12'
Set-TextValue $ws.Range("E24") '  -1.76%  '
Set-TextValue $ws.Range("E25") '  -5.09%  '
Set-TextValue $ws.Range("E26") '  -0.07%  '
Set-TextValue $ws.Range("D27") '11.30'
Set-TextValue $ws.Range("E27") '  -0.89%  '
Set-TextValue $ws.Range("D28") '3.86'
Set-TextValue $ws.Range("E28") '  -2.12%  '
Set-TextValue $ws.Range("D29") '41.17'
Set-TextValue $ws.Range("E29") '  +0.05%  '
Set-TextValue $ws.Range("E30") '  -2.26%  '
Set-TextValue $ws.Range("E31") '  -1.54%  '
Set-TextValue $ws.Range("D32") '171.58'
Set-TextValue $ws.Range("E32") '  -3.41%  '
Set-TextValue $ws.Range("D33") '21.29'
Set-TextValue $ws.Range("E33") '  -2.81%  '
Set-TextValue $ws.Range("D34") '0.0903'
Set-TextValue $ws.Range("E34") '  -1.45%  '
Set-TextValue $ws.Range("E35") '  +0.19%  '
Set-TextValue $ws.Range("E36") '  -0.45%  '
Set-TextValue $ws.Range("D37") '4.59'
Set-TextValue $ws.Range("E37") '  -4.80%  '
Set-TextValue $ws.Range("E38") '  -1.90%  '
Set-TextValue $ws.Range("D39") '3.81'
Set-TextValue $ws.Range("E39") '  -2.14%  '
Set-TextValue $ws.Range("E40") '  -8.69%  '
Set-TextValue $ws.Range("D41") '14.18'
Set-TextValue $ws.Range("E41") '  +15.97%  '
Set-TextValue $ws.Range("D42") '75.00'
Set-TextValue $ws.Range("E42") '  +11.34%  '
Set-TextValue $ws.Range("D43") '2.43'
Set-TextValue $ws.Range("E43") '  +3.23%  '
Set-TextValue $ws.Range("E44") '  -1.47%  '
Set-TextValue $ws.Range("E45") '  +9.60%  '
Set-TextValue $ws.Range("E46") '  +0.01%  '
Set-TextValue $ws.Range("E47") '  -1.78%  '
Set-TextValue $ws.Range("D48") '8.60'
Set-TextValue $ws.Range("E48") '  -2.31%  '
Set-TextValue $ws.Range("E49") '  -2.87%  '
Set-TextValue $ws.Range("D50") '100.55'
Set-TextValue $ws.Range("E50") '  +1.19%  '
Set-TextValue $ws.Range("E51") '  +0.29%  '
